# Update Name of Algo
# Apply updated numeric values to the result data produced by the
# RandomForest imputation run (Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value  = -7.949800000000003
$ws.Range("D7").Value  = -7.553199999999991
$ws.Range("C8").Value  = -11.10479999999999
$ws.Range("B12").Value = 5.543200000000001
$ws.Range("C12").Value = -14.66730000000002
$ws.Range("C14").Value = -12.5489
$ws.Range("D19").Value = -8.120699999999992
$ws.Range("E19").Value = 13.71850000000001
$ws.Range("D21").Value = -7.5356
$ws.Range("C22").Value = -10.40529999999998
$ws.Range("D24").Value = -8.311700000000002
